$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A24 was stored as a text/inline string; normalize it to a real number,
# matching what Excel does when the row is touched again.
$ws.Cells.Item(24, 1).Value = 71277628

# Append the new payment row (row 25). Columns A, B and F hold text (the
# phone number and two blank placeholder columns), so force text with a
# leading apostrophe, then reset the style so the quote-prefix marker
# doesn't linger as a cell format.
$ws.Cells.Item(25, 1).Value = "'71277628"
$ws.Cells.Item(25, 1).Style = "Normal"
$ws.Cells.Item(25, 2).Value = "'"
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 3).Value = "Cash"
$ws.Cells.Item(25, 4).Value = "2025-08-18T16:54:29"
$ws.Cells.Item(25, 5).Value = 766
$ws.Cells.Item(25, 6).Value = "'"
$ws.Cells.Item(25, 6).Style = "Normal"
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 766
